$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 311.0909
$ws.Range("I19").Value = 632.6667
$ws.Range("J19").Value = 190.5
$ws.Range("K19").Value = 632.6667
$ws.Range("L19").Value = 190.5
$ws.Range("M19").Value = -457.6667
$ws.Range("N19").Value = -540.5
$ws.Range("H131").Value = 2792.8708
$ws.Range("I131").Value = 1142.5652
$ws.Range("J131").Value = 7537.5
$ws.Range("K131").Value = 3427.6956
$ws.Range("L131").Value = 22612.5
$ws.Range("M131").Value = 1612.3044
$ws.Range("N131").Value = -32692.5
$ws.Range("H133").Value = 48035.715
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 48035.715
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 48035.715
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -58155.715
$ws.Range("H135").Value = 3617.2173
$ws.Range("I135").Value = 2340.75
$ws.Range("J135").Value = 6534.857
$ws.Range("K135").Value = 21066.75
$ws.Range("L135").Value = 58813.713
$ws.Range("M135").Value = -18531.75
$ws.Range("N135").Value = -63883.713

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 467680.9
$ws.Range("I32").Value = 5876.02
$ws.Range("J32").Value = 2391868
$ws.Range("K32").Value = 5876.02
$ws.Range("L32").Value = 2391868
$ws.Range("M32").Value = -5589.02
$ws.Range("N32").Value = -2392442
$ws.Range("H45").Value = 2353.4783
$ws.Range("I45").Value = 2106.6553
$ws.Range("J45").Value = 2774.5293
$ws.Range("K45").Value = 2106.6553
$ws.Range("L45").Value = 2774.5293
$ws.Range("M45").Value = -1729.6553
$ws.Range("N45").Value = -3528.5293
$ws.Range("H74").Value = 2365.0312
$ws.Range("I74").Value = 2186.68
$ws.Range("K74").Value = 2186.68
$ws.Range("M74").Value = -1312.68
$ws.Range("H77").Value = 2365.0312
$ws.Range("I77").Value = 2186.68
$ws.Range("K77").Value = 10933.4
$ws.Range("M77").Value = -6565.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 71133.336
$ws.Range("J140").Value = 71133.336
$ws.Range("L140").Value = 71133.336
$ws.Range("N140").Value = -81493.336

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3700.0688
$ws.Range("I31").Value = 2732.6667
$ws.Range("J31").Value = 5283.091
$ws.Range("K31").Value = 2732.6667
$ws.Range("L31").Value = 5283.091
$ws.Range("M31").Value = -2437.6667
$ws.Range("N31").Value = -5873.091
$ws.Range("H34").Value = 3700.0688
$ws.Range("I34").Value = 2732.6667
$ws.Range("J34").Value = 5283.091
$ws.Range("K34").Value = 2732.6667
$ws.Range("L34").Value = 5283.091
$ws.Range("M34").Value = -2530.6667
$ws.Range("N34").Value = -5687.091
$ws.Range("H58").Value = 125001730
$ws.Range("I58").Value = 333334140
$ws.Range("J58").Value = 2280
$ws.Range("K58").Value = 333334140
$ws.Range("L58").Value = 2280
$ws.Range("M58").Value = -333333937
$ws.Range("N58").Value = -2686
$ws.Range("H86").Value = 24159.092
$ws.Range("I86").Value = 2901.72
$ws.Range("K86").Value = 2901.72
$ws.Range("M86").Value = -1778.72
$ws.Range("H89").Value = 24159.092
$ws.Range("I89").Value = 2901.72
$ws.Range("K89").Value = 14508.6
$ws.Range("M89").Value = -8892.599999999999
$ws.Range("H136").Value = 125001730
$ws.Range("I136").Value = 333334140
$ws.Range("J136").Value = 2280
$ws.Range("K136").Value = 1000002420
$ws.Range("L136").Value = 6840
$ws.Range("M136").Value = -999999870
$ws.Range("N136").Value = -11940
$ws.Range("H140").Value = 54975
$ws.Range("I140").Value = 20000
$ws.Range("J140").Value = 89950
$ws.Range("K140").Value = 20000
$ws.Range("L140").Value = 89950
$ws.Range("M140").Value = -14820
$ws.Range("N140").Value = -100310

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1050.5
$ws.Range("I5").Value = 1325
$ws.Range("J5").Value = 867.5
$ws.Range("K5").Value = 3975
$ws.Range("L5").Value = 2602.5
$ws.Range("M5").Value = -3863
$ws.Range("N5").Value = -2826.5
$ws.Range("H122").Value = 806.6667
$ws.Range("I122").Value = 320.8
$ws.Range("K122").Value = 2887.2
$ws.Range("M122").Value = -437.2000000000003
$ws.Range("H135").Value = 1050.5
$ws.Range("I135").Value = 1325
$ws.Range("J135").Value = 867.5
$ws.Range("K135").Value = 11925
$ws.Range("L135").Value = 7807.5
$ws.Range("M135").Value = -9390
$ws.Range("N135").Value = -12877.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H122").Value = 2317.5264
$ws.Range("I122").Value = 2265.5334
$ws.Range("J122").Value = 2512.5
$ws.Range("K122").Value = 6796.600199999999
$ws.Range("L122").Value = 7537.5
$ws.Range("M122").Value = -4346.600199999999
$ws.Range("N122").Value = -12437.5
$ws.Range("H140").Value = 93929.8
$ws.Range("J140").Value = 93929.8
$ws.Range("L140").Value = 93929.8
$ws.Range("N140").Value = -104289.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3472.625
$ws.Range("I136").Value = 2151.4783
$ws.Range("J136").Value = 6848.8887
$ws.Range("K136").Value = 6454.4349
$ws.Range("L136").Value = 20546.6661
$ws.Range("M136").Value = -3904.4349
$ws.Range("N136").Value = -25646.6661

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 25416.666
$ws.Range("J112").Value = 25416.666
$ws.Range("L112").Value = 25416.666
$ws.Range("N112").Value = -28370.666
$ws.Range("H123").Value = 39597.25
$ws.Range("J123").Value = 39597.25
$ws.Range("L123").Value = 39597.25
$ws.Range("N123").Value = -49397.25
$ws.Range("H132").Value = 28849312
$ws.Range("I132").Value = 37501064
$ws.Range("J132").Value = 10141.083
$ws.Range("K132").Value = 112503192
$ws.Range("L132").Value = 30423.249
$ws.Range("M132").Value = -112500662
$ws.Range("N132").Value = -35483.249
$ws.Range("H141").Value = 70425
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360

Write-Host "All updates applied."